# Barema Seminário Sobre Microncontroladores Yduqs 31maio2023
#
# A new student (ALAN ROBERT SILVA BARROS) is added to the grade sheet as a
# new row right below row 40 (Masanori Azevedo Fukutani / "Não fez 11"
# team), and the observation note on that team is updated/expanded to
# explain that the new student delivered late and will be graded on AV3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Insert a new row 41 (pushes the old row 41 blank separator and row 42
# "Leandro de Assis Correia" entry down to rows 42/43). The inserted row
# inherits the formatting/styles of the row above it (row 40), matching
# column D (student name, style 22) and column G (observation, style 23).
$ws.Rows.Item(41).Insert()

# Update the observation text on the existing team row (row 40) and set it
# again (same text) on the newly added student's row (row 41), so both
# rows reference the same updated shared string.
$obsText = "Obs.: Aluno novo; entregou atrasado - será dado a nota na AV3(Nota da AV1)"
$ws.Cells.Item(40, 7).Value = $obsText
$ws.Cells.Item(41, 7).Value = $obsText

# New student's name in column D of the newly inserted row.
$ws.Cells.Item(41, 4).Value = "ALAN ROBERT SILVA BARROS"

# Restore the view's scroll position / selection similar to the source
# workbook (best effort - scrolled down a bit further and selecting the
# cell that now holds the "Leandro de Assis Correia" project row).
$win = $excel.ActiveWindow
$win.ScrollRow = 23
$win.ScrollColumn = 1
$ws.Range("D43").Select()
